$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the two blank placeholder rows (row 28 first so row 26's index
# is unaffected by the later deletion).
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()
